# Update the "Förändrad" (Changed) date column (C) from 2023-09-13 (45182)
# to 2023-09-15 (45184) for all data rows (2 through 92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 92 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
